$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values are swapped between the row pairs (22<->24, 23<->25)
$cols = @("A","B","E","F","G","H","Q","R","AO")

function Swap-RowValues($ws, $rowA, $rowB, $cols) {
    foreach ($col in $cols) {
        $rangeA = $ws.Range("$col$rowA")
        $rangeB = $ws.Range("$col$rowB")

        $valA = $rangeA.Value2
        $valB = $rangeB.Value2

        $rangeA.Value2 = $valB
        $rangeB.Value2 = $valA
    }
}

# Swap contents of row 22 and row 24
Swap-RowValues $ws 22 24 $cols

# Swap contents of row 23 and row 25
Swap-RowValues $ws 23 25 $cols
